# BOM-OnlyPoweredBraking.xlsx update:
# Reduce the quantity of part in row 17 (C17) from 3 to 2, which also
# recalculates the dependent formulas (D17, E-column percentages, D25 total).
# Also update the last-selected cell on Sheet1 to C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update quantity for the row-17 component; dependent formulas recalc automatically.
$ws.Range("C17").Value = 2

# Recalculate to make sure all dependent cells (D17, E2:E21, D25) are up to date.
$excel.CalculateFull()

# Match the final selection/active cell recorded in the saved sheet view.
$ws.Activate()
$ws.Range("C18").Select() | Out-Null
